# PROS-9213 - CCRU - rename KPI and add EAN code for some new Products
#
# The KPI set previously named "PoS 2019 - IC Petroleum - CAP" is renamed to
# "PoS 2019 - IC Petroleum – REG" (note: en dash, not hyphen) for the two
# "Juice Availability" atomic-name rows on the "Update Atomic Names" sheet.
# The CONCATENATE() helper formula in column E recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "PoS 2019 - IC Petroleum – REG"
$ws.Range("A23").Value = "PoS 2019 - IC Petroleum – REG"

# Reflect the author's resulting selection / scroll position.
$ws.Range("A10").Select() | Out-Null
$ws.Range("D35:D37").Select() | Out-Null
